$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert a new row (17) for the new "2508" period, pushing the footer rows down ---
$ws.Rows.Item(17).Insert()

# --- Fill the new row 17 with the same layout as row 16 (the existing "2507" entry) ---
$ws.Range("B17").Value = "CC"
$ws.Range("C17").NumberFormat = "@"
$ws.Range("C17").Value = "1047391308"
$ws.Range("D17").Value = "YEISON DE JESUS MURPHY DIAZ"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "2508"
$moneyFormat = '_-"$"\ * #,##0_-;\-"$"\ * #,##0_-;_-"$"\ * "-"??_-;_-@_-'
$ws.Range("F17").NumberFormat = $moneyFormat
$ws.Range("F17").Value = 56940
$ws.Range("G17").NumberFormat = $moneyFormat
$ws.Range("G17").Value = 1423500

# Match fonts of the new row to the row above (data row style)
$ws.Range("B17:G17").Font.Name = "Arial"
$ws.Range("H17:J17").Font.Name = "Aptos Narrow"

# Borders: thin box around every cell of the new row, like row 16
$dataRow = $ws.Range("B17:J17")
$dataRow.Borders.Item(7).LineStyle = 1
$dataRow.Borders.Item(8).LineStyle = 1
$dataRow.Borders.Item(9).LineStyle = 1
$dataRow.Borders.Item(10).LineStyle = 1
$dataRow.Borders.Item(11).LineStyle = 1

# --- Update the summary figures now that there are 2 periods in arrears ---
$ws.Range("E11").Value = 113880
$ws.Range("F13").Value = 2
